# Fix minor bugs in configuration writer
# Populate load-shift/curtailment values in columns T:X (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 0.0264
$ws.Range("U2").Value = 0.0198
$ws.Range("V2").Value = 0.0198
$ws.Range("W2").Value = 0.0066
$ws.Range("X2").Value = 0.0132
$ws.Range("T3").Value = 0.0132
$ws.Range("U3").Value = 0.0132
$ws.Range("V3").Value = 0.0132
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.0132
$ws.Range("T4").Value = 0.0132
$ws.Range("U4").Value = 0.0066
$ws.Range("V4").Value = 0.0066
$ws.Range("W4").Value = 0.0132
$ws.Range("X4").Value = 0.0132
$ws.Range("T5").Value = 0.0066
$ws.Range("U5").Value = 0.0198
$ws.Range("V5").Value = 0.0066
$ws.Range("W5").Value = 0.0066
$ws.Range("X5").Value = 0
$ws.Range("T6").Value = 0.0132
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0.0132
$ws.Range("W6").Value = 0.0066
$ws.Range("X6").Value = 0.0132
$ws.Range("T7").Value = 0.033
$ws.Range("U7").Value = 0.0132
$ws.Range("V7").Value = 0.0264
$ws.Range("W7").Value = 0.0132
$ws.Range("X7").Value = 0.0198
$ws.Range("T8").Value = 0.0198
$ws.Range("U8").Value = 0.0198
$ws.Range("V8").Value = 0.0198
$ws.Range("W8").Value = 0.0066
$ws.Range("X8").Value = 0.033
$ws.Range("T9").Value = 0.0198
$ws.Range("U9").Value = 0.0396
$ws.Range("V9").Value = 0.0066
$ws.Range("W9").Value = 0.0462
$ws.Range("X9").Value = 0.0132
$ws.Range("T10").Value = 0.0198
$ws.Range("U10").Value = 0.0264
$ws.Range("V10").Value = 0.0462
$ws.Range("W10").Value = 0.033
$ws.Range("X10").Value = 0.066
$ws.Range("T11").Value = 0.033
$ws.Range("U11").Value = 0.0726
$ws.Range("V11").Value = 0.0462
$ws.Range("W11").Value = 0.0264
$ws.Range("X11").Value = 0.033
$ws.Range("T12").Value = 0.09899999999999998
$ws.Range("U12").Value = 0.07919999999999999
$ws.Range("V12").Value = 0.05280000000000001
$ws.Range("W12").Value = 0.0462
$ws.Range("X12").Value = 0.07919999999999999
$ws.Range("T13").Value = 0.066
$ws.Range("U13").Value = 0.09239999999999998
$ws.Range("V13").Value = 0.1056
$ws.Range("W13").Value = 0.1056
$ws.Range("X13").Value = 0.07919999999999999
$ws.Range("T14").Value = 0.08579999999999999
$ws.Range("U14").Value = 0.05940000000000001
$ws.Range("V14").Value = 0.1056
$ws.Range("W14").Value = 0.09239999999999998
$ws.Range("X14").Value = 0.08579999999999999
$ws.Range("T15").Value = 0.07919999999999999
$ws.Range("U15").Value = 0.09899999999999998
$ws.Range("V15").Value = 0.07919999999999999
$ws.Range("W15").Value = 0.09239999999999998
$ws.Range("X15").Value = 0.08579999999999999
$ws.Range("T16").Value = 0.0462
$ws.Range("U16").Value = 0.1122
$ws.Range("V16").Value = 0.08579999999999999
$ws.Range("W16").Value = 0.09899999999999998
$ws.Range("X16").Value = 0.08579999999999999
$ws.Range("T17").Value = 0.1056
$ws.Range("U17").Value = 0.1056
$ws.Range("V17").Value = 0.1056
$ws.Range("W17").Value = 0.07919999999999999
$ws.Range("X17").Value = 0.1385999999999999
$ws.Range("T18").Value = 0.1385999999999999
$ws.Range("U18").Value = 0.1122
$ws.Range("V18").Value = 0.09899999999999998
$ws.Range("W18").Value = 0.09239999999999998
$ws.Range("X18").Value = 0.1188
$ws.Range("T19").Value = 0.1122
$ws.Range("U19").Value = 0.1122
$ws.Range("V19").Value = 0.132
$ws.Range("W19").Value = 0.1583999999999999
$ws.Range("X19").Value = 0.132
$ws.Range("T20").Value = 0.09239999999999998
$ws.Range("U20").Value = 0.1385999999999999
$ws.Range("V20").Value = 0.09239999999999998
$ws.Range("W20").Value = 0.1188
$ws.Range("X20").Value = 0.1517999999999999
$ws.Range("T21").Value = 0.1385999999999999
$ws.Range("U21").Value = 0.1649999999999999
$ws.Range("V21").Value = 0.07919999999999999
$ws.Range("W21").Value = 0.1188
$ws.Range("X21").Value = 0.1385999999999999
$ws.Range("T22").Value = 0.1122
$ws.Range("U22").Value = 0.1517999999999999
$ws.Range("V22").Value = 0.0726
$ws.Range("W22").Value = 0.09239999999999998
$ws.Range("X22").Value = 0.066
$ws.Range("T23").Value = 0.09239999999999998
$ws.Range("U23").Value = 0.09899999999999998
$ws.Range("V23").Value = 0.05280000000000001
$ws.Range("W23").Value = 0.0726
$ws.Range("X23").Value = 0.0726
$ws.Range("T24").Value = 0.0462
$ws.Range("U24").Value = 0.066
$ws.Range("V24").Value = 0.05280000000000001
$ws.Range("W24").Value = 0.0198
$ws.Range("X24").Value = 0.0726
$ws.Range("T25").Value = 0.0198
$ws.Range("U25").Value = 0.033
$ws.Range("V25").Value = 0.0264
$ws.Range("W25").Value = 0.0198
$ws.Range("X25").Value = 0.033
